$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-12-21 Saturday" "2024-12-22 Sunday"

Replace-Text "50÷8=" "12÷8="
Replace-Text "48÷8=" "73÷4="
Replace-Text "88÷8=" "79÷9="
Replace-Text "94÷3=" "26÷6="
Replace-Text "61÷9=" "47÷6="

Replace-Text "77÷4=" "32÷2="
Replace-Text "10÷8=" "69÷7="
Replace-Text "32÷8=" "53÷9="
Replace-Text "37÷4=" "45÷5="
Replace-Text "17÷2=" "35÷4="

Replace-Text "12÷9=" "52÷6="
Replace-Text "25÷8=" "23÷5="
Replace-Text "62÷9=" "32÷7="
Replace-Text "66÷3=" "25÷8="
Replace-Text "77÷9=" "37÷4="

Replace-Text "95÷8=" "62÷5="
Replace-Text "10÷3=" "63÷8="
Replace-Text "75÷2=" "37÷9="
Replace-Text "60÷5=" "50÷9="
Replace-Text "28÷9=" "10÷2="

Replace-Text "14÷3=" "41÷2="
Replace-Text "79÷3=" "19÷9="
Replace-Text "26÷5=" "40÷4="
Replace-Text "47÷3=" "98÷6="
Replace-Text "19÷7=" "77÷4="
